# Refresh the cryptocurrency price/volume snapshot on Sheet1 (rows 2-51).
# Numeric-looking strings in column D (Price) must stay TEXT, matching the
# source data's formatting (e.g. thousands separators like '61.455.17'), so
# we briefly force a Text number format before assigning, then restore the
# default 'Normal' style so no visible formatting changes remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

# Row 2
$ws.Range('D2').Value = '61.455.17'
$ws.Range('E2').Value = '  -4.38%  '
# Row 3
$ws.Range('D3').Value = '2.971.06'
$ws.Range('E3').Value = '  -5.28%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
Set-TextValue 'D5' '538.24'
$ws.Range('E5').Value = '  -5.86%  '
# Row 6
Set-TextValue 'D6' '149.22'
$ws.Range('E6').Value = '  -8.74%  '
# Row 7
$ws.Range('E7').Value = '  +0.02%  '
# Row 8
Set-TextValue 'D8' '0.567'
$ws.Range('E8').Value = '  -1.11%  '
# Row 9
$ws.Range('D9').Value = '2.980.00'
$ws.Range('E9').Value = '  -5.41%  '
# Row 10
$ws.Range('E10').Value = '  -3.72%  '
# Row 11
Set-TextValue 'D11' '6.13'
$ws.Range('E11').Value = '  -7.10%  '
# Row 12
$ws.Range('E12').Value = '  -4.78%  '
# Row 13
$ws.Range('D13').Value = '3.486.77'
$ws.Range('E13').Value = '  -5.38%  '
# Row 14
$ws.Range('E14').Value = '  -1.62%  '
# Row 15
$ws.Range('D15').Value = '61.535.07'
$ws.Range('E15').Value = '  -4.29%  '
# Row 16
Set-TextValue 'D16' '23.57'
$ws.Range('E16').Value = '  -6.05%  '
# Row 17
$ws.Range('D17').Value = '2.973.24'
$ws.Range('E17').Value = '  -5.48%  '
# Row 18
$ws.Range('E18').Value = '  -5.43%  '
# Row 19
Set-TextValue 'D19' '5.15'
$ws.Range('E19').Value = '  -1.93%  '
# Row 20
Set-TextValue 'D20' '12.04'
$ws.Range('E20').Value = '  -3.81%  '
# Row 21
Set-TextValue 'D21' '378.57'
$ws.Range('E21').Value = '  -5.76%  '
# Row 22
Set-TextValue 'D22' '6.68'
$ws.Range('E22').Value = '  -5.53%  '
# Row 23
$ws.Range('E23').Value = '  +0.22%  '
# Row 24
Set-TextValue 'D24' '5.65'
$ws.Range('E24').Value = '  -3.68%  '
# Row 25
Set-TextValue 'D25' '65.56'
$ws.Range('E25').Value = '  -4.55%  '
# Row 26
Set-TextValue 'D26' '0.470'
# Row 27
$ws.Range('D27').Value = '3.093.11'
$ws.Range('E27').Value = '  -5.63%  '
# Row 28
$ws.Range('E28').Value = '  -5.05%  '
# Row 29
$ws.Range('E29').Value = '  +0.20%  '
# Row 30
$ws.Range('D30').Value = '0.0₃0935'
$ws.Range('E30').Value = '  -7.42%  '
# Row 31
Set-TextValue 'D31' '8.19'
$ws.Range('E31').Value = '  -6.61%  '
# Row 32
$ws.Range('E32').Value = '  +0.01%  '
# Row 33
$ws.Range('E33').Value = '  -5.36%  '
# Row 34
Set-TextValue 'D34' '20.39'
$ws.Range('E34').Value = '  -3.98%  '
# Row 35
Set-TextValue 'D35' '159.53'
$ws.Range('E35').Value = '  -1.26%  '
# Row 36
$ws.Range('E36').Value = '  -4.53%  '
# Row 37
Set-TextValue 'D37' '5.91'
$ws.Range('E37').Value = '  -5.90%  '
# Row 38
$ws.Range('E38').Value = '  -3.86%  '
# Row 39
$ws.Range('E39').Value = '  -5.75%  '
# Row 40
$ws.Range('E40').Value = '  -7.89%  '
# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D41' '3.91'
$ws.Range('E41').Value = '  -3.83%  '
# Row 42
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D42' '37.49'
$ws.Range('E42').Value = '  -2.27%  '
# Row 43
$ws.Range('D43').Value = '2.412.52'
$ws.Range('E43').Value = '  -8.57%  '
# Row 44
$ws.Range('E44').Value = '  -6.93%  '
# Row 45
Set-TextValue 'D45' '0.669'
$ws.Range('E45').Value = '  -2.93%  '
# Row 46
Set-TextValue 'D46' '0.0590'
$ws.Range('E46').Value = '  -3.90%  '
# Row 47
$ws.Range('E47').Value = '  +0.02%  '
# Row 48
Set-TextValue 'D48' '5.03'
$ws.Range('E48').Value = '  -7.25%  '
# Row 49
Set-TextValue 'D49' '0.0245'
$ws.Range('E49').Value = '  -3.64%  '
# Row 50
Set-TextValue 'D50' '0.0950'
$ws.Range('E50').Value = '  -2.54%  '
# Row 51
Set-TextValue 'D51' '19.71'
$ws.Range('E51').Value = '  -6.76%  '
